# Sprint 3 planning sheet update — product backlog v5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Contato Usuário" story block (rows 11-17) becomes "Contatar Profissional" ---
$ws.Range("A11").Value = "Estória de Usuário: CH:2 ID:11 – Contatar Profissional"

# --- "Mural de Serviços" story block (rows 1-9) becomes "Publicar no Mural" ---
$ws.Range("A1").Value = "Estória de Usuário: CH:2 ID:11 – Publicar no Mural"
$ws.Range("B4").Value = "Criar tela publicarNoMural"
$ws.Range("G4").Value = "Com campos obrigatórios título, descricao e categoria que será uma lista das categoria do sistema. E o campo contato não obrigatório"
$ws.Range("G5").Value = "Com método publicarNoMural passando por parametro o id da categoria selecionada."
$ws.Range("G7").Value = "Com método publicarNoMural passando por parametro o id da categoria selecionada. "

# The longer wrapped text in G4 makes row 4 auto-fit taller
$ws.Rows.Item(4).AutoFit()

# --- View state: scroll back to top, move selection to F12 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F12").Select()
